# ============================================================
# Edit: add "2022-Q3" sheet (fund holdings) right after "总计",
# shifting 2022-Q2..2020-Q4 down one slot each; update the "总计"
# summary sheet with the new quarter row.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- Part 1: insert the new worksheet "2022-Q3" right after "总计" ----
$zongji = $wb.Worksheets.Item("总计")
$q2old  = $wb.Worksheets.Item("2022-Q2")   # existing fund sheet we borrow formatting from
$newWs  = $wb.Worksheets.Add($null, $zongji)
$newWs.Name = "2022-Q3"

# Borrow the header-row (B1:H1) and index-column (A) styling (bold + border,
# style index "2" in this workbook) by copying from the sheet that already
# uses that exact layout, so the new sheet reuses the same style record
# instead of manufacturing a new (visually-equivalent) one.
$q2old.Range("A1:H21").Copy($newWs.Range("A1:H21"))

$headers = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
$rows = @(
  ,@('008969', '睿远均衡价值三年持有期混合A', '121.89', '90.22', '3.64', '4.4368', 9)
  ,@('000696', '汇添富环保行业股票', '41.97', '91.56', '3.15', '1.3221', 10)
  ,@('014362', '睿远稳进配置两年持有混合A', '64.40', '35.09', '1.33', '0.8565', 8)
  ,@('008970', '睿远均衡价值三年持有期混合C', '15.10', '90.22', '3.64', '0.5496', 9)
  ,@('011296', '汇添富优势行业一年定开混合A', '11.39', '86.24', '4.67', '0.5319', 6)
  ,@('014363', '睿远稳进配置两年持有混合C', '35.67', '35.09', '1.33', '0.4744', 8)
  ,@('011122', '汇添富ESG可持续成长股票A', '9.39', '82.20', '4.76', '0.4470', 6)
  ,@('517160', '南方中证长江保护主题ETF', '16.84', '99.30', '2.55', '0.4294', 5)
  ,@('517330', '易方达中证长江保护主题ETF', '16.37', '99.46', '2.56', '0.4191', 5)
  ,@('001166', '建信环保产业股票', '7.82', '88.64', '4.28', '0.3347', 7)
  ,@('013365', '汇添富产业升级混合A', '4.34', '84.74', '2.53', '0.1098', 10)
  ,@('501030', '汇添富中证环境治理指数（LOF）A', '3.12', '92.74', '2.22', '0.0693', 5)
  ,@('014522', '汇添富低碳投资一年持有混合A', '2.07', '85.02', '2.78', '0.0575', 7)
  ,@('001742', '广发百发大数据策略精选灵活配置混合E', '2.51', '40.85', '2.12', '0.0532', 10)
  ,@('164908', '交银施罗德中证环境治理指数（LOF）', '1.57', '93.62', '2.20', '0.0345', 5)
  ,@('002634', '华宝未来主导产业灵活配置混合A', '0.46', '91.91', '7.08', '0.0326', 4)
  ,@('501031', '汇添富中证环境治理指数（LOF）C', '1.30', '92.74', '2.22', '0.0289', 5)
  ,@('011123', '汇添富ESG可持续成长股票C', '0.44', '82.20', '4.76', '0.0209', 6)
  ,@('164401', '前海开源中证健康产业指数', '1.81', '94.19', '1.13', '0.0205', 8)
  ,@('014523', '汇添富低碳投资一年持有混合C', '0.51', '85.02', '2.78', '0.0142', 7)
  ,@('013366', '汇添富产业升级混合C', '0.29', '84.74', '2.53', '0.0073', 10)
  ,@('011297', '汇添富优势行业一年定开混合C', '0.15', '86.24', '4.67', '0.0070', 6)
  ,@('012919', '华宝未来主导产业灵活配置混合C', '0.07', '91.91', '7.08', '0.0050', 4)
  ,@('001741', '广发百发大数据策略精选灵活配置混合A', '0.21', '40.85', '2.12', '0.0045', 10)
  ,@('013413', '交银施罗德中证环境治理指数（LOF）C', '0.09', '93.62', '2.20', '0.0020', 5)
)

# Header row text (overwrite the copied text with itself / ensure labels match)
for ($col = 2; $col -le 8; $col++) {
    $newWs.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Force columns B..G to store values as literal text (matches the source
# workbook, which keeps these as text, not numbers) - column A (index) and
# H (rank) stay numeric.
$dataRange = $newWs.Range("B2:G26")
$dataRange.NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $newWs.Cells.Item($r, 1).Value = $r - 2          # 0-based running index in col A
    $newWs.Cells.Item($r, 2).Value = $row[0]
    $newWs.Cells.Item($r, 3).Value = $row[1]
    $newWs.Cells.Item($r, 4).Value = $row[2]
    $newWs.Cells.Item($r, 5).Value = $row[3]
    $newWs.Cells.Item($r, 6).Value = $row[4]
    $newWs.Cells.Item($r, 7).Value = $row[5]
    $newWs.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# The donor sheet (2022-Q2) only had 20 data rows (21 used rows incl.
# header); our new sheet has 25 data rows (26 used rows), so clear any
# leftover copied formatting/content beyond row 21 isn't an issue since we
# copied exactly A1:H21 - rows 22..26 still need the plain/no-style layout
# that the rest of the data rows use (B:H un-styled, A styled "2").
$q2old.Range("A21").Copy($newWs.Range("A22:A26"))

# ---- Part 2: update the "总计" (summary) sheet ----
$summary = $wb.Worksheets.Item("总计")

# Create the new row 9 by copying row 8's formatting (keeps column-A's bold
# border style) then shifting all the quarterly rows down by one and
# writing the new 2022-Q3 figures into row 2 (row 2's own cells/style are
# left as-is other than the text/number updates).
$summary.Range("A8").Copy($summary.Range("A9"))
$summary.Cells.Item(9, 1).Value = 7

for ($row = 8; $row -ge 2; $row--) {
    $dst = $row + 1
    $summary.Cells.Item($dst, 2).Value = $summary.Cells.Item($row, 2).Value()
    $summary.Cells.Item($dst, 3).Value = $summary.Cells.Item($row, 3).Value()
    $summary.Cells.Item($dst, 4).Value = $summary.Cells.Item($row, 4).Value()
}

$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 25
$summary.Cells.Item(2, 4).Value = 10.27
